$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.247929930686951
$ws.Range("B1").Value = 2.386881351470947
$ws.Range("C1").Value = 3.991010189056396
$ws.Range("D1").Value = 2.724275588989258
$ws.Range("E1").Value = 1.317806005477905
